# C5-PowerPoint.pptx edit — Sun, May 24, 2020 12:05:19 AM
#
# The source diff shows the table on slide 6 (the "SOURCES OF FINANCE"
# slide) being re-styled: its custom table style
# {4F7B065C-77D5-40D3-B897-1C6C75F0044C} ("Table_0") is swapped for the
# built-in PowerPoint table style {40030BFF-4631-4D82-82C5-25584FF95ABB}.
#
# (The diff also shows the deck's embedded-font declaration disappearing
# from ppt/presentation.xml — that is a side effect of how PowerPoint
# persists the "Embed fonts in the file" Save option and isn't something
# exposed anywhere on the Presentation/Application COM object model, so
# it can't be driven from automation code; it is not reproduced here.)

$p = $ppt.ActivePresentation

# Slide 6 holds a single table ("Google Shape;127;p18") as its second shape.
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)

if (-not $tableShape.HasTable) {
    throw "Expected shape 2 on slide 6 to contain a table"
}

$table = $tableShape.Table
$table.ApplyStyle("{40030BFF-4631-4D82-82C5-25584FF95ABB}")
